# LM3150MHX Calculator - rename 17V boards to 18V.
#
# The "12V" tab (3rd sheet, sheetId 3) models a board whose typical input
# voltage (B9, "Vintyp") was mislabeled as 17V; it should read 18V. All the
# other numbers that change in this workbook are formulas that depend on
# B9 (directly or indirectly), so Excel's automatic recalculation takes
# care of them once B9 is corrected here - no other cell needs to be
# touched.
#
# The commit also nudges the saved "current view" of the workbook: the
# 12V sheet becomes the active/selected tab (with B10 selected) instead
# of the 5V sheet, which matches the author simply having been looking
# at the 12V tab (cell B10) when they saved the file after the edit.

$wb = $excel.ActiveWorkbook

$ws5V   = $wb.Worksheets.Item(1)   # "5V"
$ws33V  = $wb.Worksheets.Item(2)   # "3.3V"
$ws12V  = $wb.Worksheets.Item(3)   # "12V"

# --- the actual data fix: 17V -> 18V -------------------------------------
$ws12V.Range("B9").Value = 18

# --- saved view state: make the 12V sheet the active tab, with B10 -------
# selected (mirrors the author's on-screen state when the file was saved;
# activating it last also makes it the workbook's active/selected tab).
$ws12V.Activate()
$ws12V.Range("B10").Select()
